# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Column G ("K") values for rows 2-29 are recalculated/rewritten with the
# correct strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 7
    4  = 5
    5  = 6
    6  = 3
    7  = 4
    8  = 11
    9  = 6
    10 = 3
    11 = 2
    12 = 3
    13 = 5
    14 = 5
    15 = 6
    16 = 4
    17 = 5
    18 = 2
    19 = 7
    20 = 7
    21 = 5
    22 = 2
    23 = 8
    24 = 4
    25 = 5
    26 = 6
    27 = 5
    28 = 2
    29 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
